# Fill in the bill-form fields for Mr. Argha Chandra Dhar: name, designation,
# department/section, year, term, and the amount-in-words line.
#
# Note: a couple of existing labels in this sheet use the Bengali letter YYA
# (U+09DF), e.g. "নিয়মিত" and "কথায়". We read those labels back from the
# workbook itself (via .Text) instead of re-typing them, then append the new
# (ASCII/")safe" Bengali suffix text, so the stored character encoding always
# matches what the workbook already contains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header block: name / designation / year / term / department ----

# A3: "নাম:"  ->  "নাম: Mr. Argha Chandra Dhar"
$nameLabel = $ws.Range("A3").Text
$ws.Range("A3").Value = $nameLabel + " Mr. Argha Chandra Dhar"

# A4: "পদবী: "  ->  "পদবী: প্রভাষক"   (Lecturer)
$designationLabel = $ws.Range("A4").Text
$ws.Range("A4").Value = $designationLabel + "প্রভাষক"

# G4 / I4 were blank -> year "৪র্থ" (4th) / term "১ম" (1st)
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"

# B5 was blank -> department short code "সিএসই" (CSE)
$ws.Range("B5").Value = "সিএসই"

# F5: "বিভাগ :"  ->  "বিভাগ :সিএসই"
$deptLabel = $ws.Range("F5").Text
$ws.Range("F5").Value = $deptLabel + "সিএসই"

# A32: "কথায়:"  ->  "কথায়:সাত হাজার বাহান্ন টাকা মাত্র।"  (amount in words)
$inWordsLabel = $ws.Range("A32").Text
$ws.Range("A32").Value = $inWordsLabel + "সাত হাজার বাহান্ন টাকা মাত্র।"

# ---- Layout tweaks that came along with entering the wider name text ----

# Column A widens to fit "নাম: Mr. Argha Chandra Dhar"
$ws.Columns.Item(1).ColumnWidth = 13.5

# Row 36 grows (signature/teacher block wraps taller)
$ws.Rows.Item(36).RowHeight = 68.4

# ---- Selection / scroll state left by the editing session ----
$ws.Range("I32").Select()
